# TC-082: add a new "tc082" worksheet (authored test-case data) right after
# "tc076" and before "tc048", matching the AddTest/tc0xx sheet template used
# throughout this workbook. The new sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("tc076")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "tc082"

# Header row
$newSheet.Range("A1").Value = "Epic"
$newSheet.Range("B1").Value = "Feature"
$newSheet.Range("C1").Value = "Requirement"
$newSheet.Range("D1").Value = "Tcname"
$newSheet.Range("E1").Value = "desc"
$newSheet.Range("F1").Value = "Type"
$newSheet.Range("G1").Value = "assigned"

# Data row - note: E2 is written before D2 so the two brand-new shared
# strings land in the same table order as the authored workbook.
$newSheet.Range("A2").Value = "Epic Mohit"
$newSheet.Range("B2").Value = "Mohit Feature"
$newSheet.Range("C2").Value = "RQ-489"
$newSheet.Range("E2").Value = "Creating testcase for automation TC-082"
$newSheet.Range("D2").Value = "Testing Notification for create TC Jan 05-01-2026"
$newSheet.Range("F2").Value = "Manual"
$newSheet.Range("G2").Value = "Mohit Aman"

# Column widths for the (now much longer) Tcname/desc columns
$newSheet.Columns.Item(4).ColumnWidth = 38.5
$newSheet.Columns.Item(5).ColumnWidth = 36.333333333333336

# Make the new sheet the active tab with A1:G2 selected
$newSheet.Range("A1:G2").Select() | Out-Null
